$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 131, shifting existing rows 131:193 down to 132:194.
$ws.Rows.Item(131).Insert()

# Populate the newly inserted row 131 with the new "Femacal de La Calera" /
# Zanahoria record (values mirror the surrounding rows' constant columns).
$ws.Range("A131").Value = 3
$ws.Range("B131").Value = "Femacal de La Calera"
$ws.Range("C131").Value = "Coquimbo"
$ws.Range("D131").Value = 44466
$ws.Range("E131").Value = 5
$ws.Range("F131").Value = 100114013
$ws.Range("G131").Value = "Zanahoria"
$ws.Range("H131").Value = "Sin especificar"
$ws.Range("I131").Value = "Primera"
$ws.Range("J131").Value = 370
$ws.Range("K131").Value = 6500
$ws.Range("L131").Value = 7000
$ws.Range("M131").Value = 6743
$ws.Range("N131").Value = "$/saco 20 kilos"
$ws.Range("O131").Value = "Provincia de Quillota"
$ws.Range("P131").Value = 337
$ws.Range("Q131").Value = 20
$ws.Range("R131").Value = "Hortaliza"
